$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 53.64296733333333
$ws.Range("H2").Value = 160.928902
$ws.Range("I2").Value = 0.1405570576660657
$ws.Range("J2").Value = 0.1405570576660657
$ws.Range("M2").Value = 0.4890553333333333
$ws.Range("N2").Value = 1.467166
$ws.Range("O2").Value = 0.9644476581758422
$ws.Range("P2").Value = 0.9644476581758422
$ws.Range("Q2").Value = 26.23437927019244
$ws.Range("R2").Value = 236.109413431732
$ws.Range("S2").Value = 0.1355599251061239
$ws.Range("T2").Value = 0.1355599251061239

# Row 3
$ws.Range("G3").Value = 53.64296733333333
$ws.Range("H3").Value = 160.928902
$ws.Range("I3").Value = 0.1405570576660657
$ws.Range("J3").Value = 0.1405570576660657
$ws.Range("O3").Value = 0.03555234182415776
$ws.Range("P3").Value = 0.03555234182415776
$ws.Range("Q3").Value = 0.9670754150853332
$ws.Range("R3").Value = 8.703678735767999
$ws.Range("S3").Value = 0.004997132559941823
$ws.Range("T3").Value = 0.004997132559941823

# Row 4
$ws.Range("I4").Value = 0.83973167405618
$ws.Range("J4").Value = 0.8397316740561799
$ws.Range("M4").Value = 0.4890553333333333
$ws.Range("N4").Value = 1.467166
$ws.Range("O4").Value = 0.9644476581758422
$ws.Range("P4").Value = 0.9644476581758422
$ws.Range("Q4").Value = 156.7323590020058
$ws.Range("R4").Value = 1410.591231018052
$ws.Range("S4").Value = 0.8098772465395624
$ws.Range("T4").Value = 0.8098772465395623

# Row 5
$ws.Range("I5").Value = 0.83973167405618
$ws.Range("J5").Value = 0.8397316740561799
$ws.Range("O5").Value = 0.03555234182415776
$ws.Range("P5").Value = 0.03555234182415776
$ws.Range("S5").Value = 0.02985442751661754
$ws.Range("T5").Value = 0.02985442751661754

# Row 6
$ws.Range("G6").Value = 7.522716666666668
$ws.Range("I6").Value = 0.01971126827775425
$ws.Range("J6").Value = 0.01971126827775425
$ws.Range("M6").Value = 0.4890553333333333
$ws.Range("N6").Value = 1.467166
$ws.Range("O6").Value = 0.9644476581758422
$ws.Range("P6").Value = 0.9644476581758422
$ws.Range("Q6").Value = 3.679024706988889
$ws.Range("R6").Value = 33.1112223629
$ws.Range("S6").Value = 0.01901048653015585
$ws.Range("T6").Value = 0.01901048653015585

# Row 7
$ws.Range("G7").Value = 7.522716666666668
$ws.Range("I7").Value = 0.01971126827775425
$ws.Range("J7").Value = 0.01971126827775425
$ws.Range("O7").Value = 0.03555234182415776
$ws.Range("P7").Value = 0.03555234182415776
$ws.Range("S7").Value = 0.0007007817475983966
$ws.Range("T7").Value = 0.0007007817475983965
